# Update "想去人数" (want-to-go count) and, where applicable, "最低票价"
# (lowest ticket price) figures across the four sheets of the workbook,
# reflecting a refreshed data pull (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8167
$ws1.Range("F4").Value = 1919
$ws1.Range("G4").Value = 68
$ws1.Range("F5").Value = 6516
$ws1.Range("F7").Value = 2074
$ws1.Range("F8").Value = 572
$ws1.Range("F9").Value = 48
$ws1.Range("F10").Value = 20
$ws1.Range("F11").Value = 50
$ws1.Range("F14").Value = 65
$ws1.Range("F15").Value = 3
$ws1.Range("F16").Value = 8543
$ws1.Range("F19").Value = 194
$ws1.Range("F21").Value = 1808
$ws1.Range("F26").Value = 37
$ws1.Range("F28").Value = 186
$ws1.Range("F29").Value = 591
$ws1.Range("F31").Value = 16
$ws1.Range("F32").Value = 5
$ws1.Range("F33").Value = 2088
$ws1.Range("F34").Value = 848
$ws1.Range("F35").Value = 474
$ws1.Range("F39").Value = 185
$ws1.Range("F40").Value = 148
$ws1.Range("F42").Value = 28
$ws1.Range("F44").Value = 31

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 391
$ws2.Range("F7").Value = 441
$ws2.Range("F16").Value = 97
$ws2.Range("F22").Value = 66

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 2323
$ws3.Range("F4").Value = 308

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2323
$ws4.Range("F4").Value = 391
$ws4.Range("F5").Value = 8167
$ws4.Range("F7").Value = 308
$ws4.Range("F8").Value = 1919
$ws4.Range("G8").Value = 68
$ws4.Range("F9").Value = 6516
$ws4.Range("F11").Value = 2074
$ws4.Range("F14").Value = 572
$ws4.Range("F15").Value = 48
$ws4.Range("F18").Value = 50
$ws4.Range("F22").Value = 65
$ws4.Range("F23").Value = 8543
$ws4.Range("F25").Value = 194
$ws4.Range("F27").Value = 1808
$ws4.Range("F30").Value = 186
$ws4.Range("F31").Value = 591
$ws4.Range("F32").Value = 2088
$ws4.Range("F33").Value = 848
$ws4.Range("F35").Value = 474
$ws4.Range("F40").Value = 97
$ws4.Range("F45").Value = 66
